$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 48.484375
$ws.Range("H2").Value = 145.453125
$ws.Range("I2").Value = 0.7776469276297807
$ws.Range("J2").Value = 0.7776469276297806
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 0.3363556666666667
$ws.Range("N2").Value = 1.009067
$ws.Range("O2").Value = 0.01004222540086594
$ws.Range("P2").Value = 0.01004222540086594
$ws.Range("Q2").Value = 16.30799427604166
$ws.Range("R2").Value = 146.771948484375
$ws.Range("S2").Value = 0.00780930572954914
$ws.Range("T2").Value = 0.007809305729549137

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 48.484375
$ws.Range("H3").Value = 145.453125
$ws.Range("I3").Value = 0.7776469276297807
$ws.Range("J3").Value = 0.7776469276297806
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 5.511159333333334
$ws.Range("N3").Value = 16.533478
$ws.Range("O3").Value = 0.1645410193141369
$ws.Range("P3").Value = 0.1645410193141369
$ws.Range("Q3").Value = 267.2051158020834
$ws.Range("R3").Value = 2404.84604221875
$ws.Range("S3").Value = 0.127954818138711
$ws.Range("T3").Value = 0.127954818138711

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 48.484375
$ws.Range("H4").Value = 145.453125
$ws.Range("I4").Value = 0.7776469276297807
$ws.Range("J4").Value = 0.7776469276297806
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 26.26335166666667
$ws.Range("N4").Value = 78.790055
$ws.Range("O4").Value = 0.784117894705331
$ws.Range("P4").Value = 0.7841178947053309
$ws.Range("Q4").Value = 1273.362190963541
$ws.Range("R4").Value = 11460.25971867187
$ws.Range("S4").Value = 0.6097668717171325
$ws.Range("T4").Value = 0.6097668717171324

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 48.484375
$ws.Range("H5").Value = 145.453125
$ws.Range("I5").Value = 0.7776469276297807
$ws.Range("J5").Value = 0.7776469276297806
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 1.383269666666667
$ws.Range("N5").Value = 4.149809
$ws.Range("O5").Value = 0.04129886057966625
$ws.Range("P5").Value = 0.04129886057966624
$ws.Range("Q5").Value = 67.06696524479167
$ws.Range("R5").Value = 603.602687203125
$ws.Range("S5").Value = 0.03211593204438812
$ws.Range("T5").Value = 0.03211593204438811

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 11.47148733333333
$ws.Range("H6").Value = 34.414462
$ws.Range("I6").Value = 0.1839926137051496
$ws.Range("J6").Value = 0.1839926137051496
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 0.3363556666666667
$ws.Range("N6").Value = 1.009067
$ws.Range("O6").Value = 0.01004222540086594
$ws.Range("P6").Value = 0.01004222540086594
$ws.Range("Q6").Value = 3.858499769661556
$ws.Range("R6").Value = 34.726497926954
$ws.Range("S6").Value = 0.001847695298921568
$ws.Range("T6").Value = 0.001847695298921567

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 11.47148733333333
$ws.Range("H7").Value = 34.414462
$ws.Range("I7").Value = 0.1839926137051496
$ws.Range("J7").Value = 0.1839926137051496
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 5.511159333333334
$ws.Range("N7").Value = 16.533478
$ws.Range("O7").Value = 0.1645410193141369
$ws.Range("P7").Value = 0.1645410193141369
$ws.Range("Q7").Value = 63.22119448431513
$ws.Range("R7").Value = 568.9907503588361
$ws.Range("S7").Value = 0.03027433220531756
$ws.Range("T7").Value = 0.03027433220531755

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 11.47148733333333
$ws.Range("H8").Value = 34.414462
$ws.Range("I8").Value = 0.1839926137051496
$ws.Range("J8").Value = 0.1839926137051496
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 26.26335166666667
$ws.Range("N8").Value = 78.790055
$ws.Range("O8").Value = 0.784117894705331
$ws.Range("P8").Value = 0.7841178947053309
$ws.Range("Q8").Value = 301.2797059750455
$ws.Range("R8").Value = 2711.51735377541
$ws.Range("S8").Value = 0.1442719008998131
$ws.Range("T8").Value = 0.1442719008998131

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 11.47148733333333
$ws.Range("H9").Value = 34.414462
$ws.Range("I9").Value = 0.1839926137051496
$ws.Range("J9").Value = 0.1839926137051496
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 1.383269666666667
$ws.Range("N9").Value = 4.149809
$ws.Range("O9").Value = 0.04129886057966625
$ws.Range("P9").Value = 0.04129886057966624
$ws.Range("Q9").Value = 15.86816045975089
$ws.Range("R9").Value = 142.813444137758
$ws.Range("S9").Value = 0.007598685301097362
$ws.Range("T9").Value = 0.00759868530109736

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 0.356432
$ws.Range("H10").Value = 1.069296
$ws.Range("I10").Value = 0.005716857228930723
$ws.Range("J10").Value = 0.005716857228930722
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 0.3363556666666667
$ws.Range("N10").Value = 1.009067
$ws.Range("O10").Value = 0.01004222540086594
$ws.Range("P10").Value = 0.01004222540086594
$ws.Range("Q10").Value = 0.1198879229813333
$ws.Range("R10").Value = 1.078991306832
$ws.Range("S10").Value = 0.00005740996887749216
$ws.Range("T10").Value = 0.00005740996887749214

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 0.356432
$ws.Range("H11").Value = 1.069296
$ws.Range("I11").Value = 0.005716857228930723
$ws.Range("J11").Value = 0.005716857228930722
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 5.511159333333334
$ws.Range("N11").Value = 16.533478
$ws.Range("O11").Value = 0.1645410193141369
$ws.Range("P11").Value = 0.1645410193141369
$ws.Range("Q11").Value = 1.964353543498667
$ws.Range("R11").Value = 17.679181891488
$ws.Range("S11").Value = 0.0009406575157216533
$ws.Range("T11").Value = 0.0009406575157216531

$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 0.356432
$ws.Range("H12").Value = 1.069296
$ws.Range("I12").Value = 0.005716857228930723
$ws.Range("J12").Value = 0.005716857228930722
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 26.26335166666667
$ws.Range("N12").Value = 78.790055
$ws.Range("O12").Value = 0.784117894705331
$ws.Range("P12").Value = 0.7841178947053309
$ws.Range("Q12").Value = 9.361098961253333
$ws.Range("R12").Value = 84.24989065128
$ws.Range("S12").Value = 0.004482690054680111
$ws.Range("T12").Value = 0.00448269005468011

$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 0.356432
$ws.Range("H13").Value = 1.069296
$ws.Range("I13").Value = 0.005716857228930723
$ws.Range("J13").Value = 0.005716857228930722
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 1.383269666666667
$ws.Range("N13").Value = 4.149809
$ws.Range("O13").Value = 0.04129886057966625
$ws.Range("P13").Value = 0.04129886057966624
$ws.Range("Q13").Value = 0.4930415738293334
$ws.Range("R13").Value = 4.437374164464001
$ws.Range("S13").Value = 0.000236099689651467
$ws.Range("T13").Value = 0.000236099689651467

$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 2.035248333333333
$ws.Range("H14").Value = 6.105745
$ws.Range("I14").Value = 0.03264360143613892
$ws.Range("J14").Value = 0.03264360143613892
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 0.3363556666666667
$ws.Range("N14").Value = 1.009067
$ws.Range("O14").Value = 0.01004222540086594
$ws.Range("P14").Value = 0.01004222540086594
$ws.Range("Q14").Value = 0.6845673099905555
$ws.Range("R14").Value = 6.161105789914999
$ws.Range("S14").Value = 0.0003278144035177381
$ws.Range("T14").Value = 0.0003278144035177381

$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 2.035248333333333
$ws.Range("H15").Value = 6.105745
$ws.Range("I15").Value = 0.03264360143613892
$ws.Range("J15").Value = 0.03264360143613892
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 5.511159333333334
$ws.Range("N15").Value = 16.533478
$ws.Range("O15").Value = 0.1645410193141369
$ws.Range("P15").Value = 0.1645410193141369
$ws.Range("Q15").Value = 11.21657784790111
$ws.Range("R15").Value = 100.94920063111
$ws.Range("S15").Value = 0.005371211454386723
$ws.Range("T15").Value = 0.005371211454386722

$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 2.035248333333333
$ws.Range("H16").Value = 6.105745
$ws.Range("I16").Value = 0.03264360143613892
$ws.Range("J16").Value = 0.03264360143613892
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 26.26335166666667
$ws.Range("N16").Value = 78.790055
$ws.Range("O16").Value = 0.784117894705331
$ws.Range("P16").Value = 0.7841178947053309
$ws.Range("Q16").Value = 53.45244270733055
$ws.Range("R16").Value = 481.071984365975
$ws.Range("S16").Value = 0.02559643203370517
$ws.Range("T16").Value = 0.02559643203370517

$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 2.035248333333333
$ws.Range("H17").Value = 6.105745
$ws.Range("I17").Value = 0.03264360143613892
$ws.Range("J17").Value = 0.03264360143613892
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 1.383269666666667
$ws.Range("N17").Value = 4.149809
$ws.Range("O17").Value = 0.04129886057966625
$ws.Range("P17").Value = 0.04129886057966624
$ws.Range("Q17").Value = 2.815297283633889
$ws.Range("R17").Value = 25.337675552705
$ws.Range("S17").Value = 0.001348143544529294
$ws.Range("T17").Value = 0.001348143544529294
